$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, pushing the existing row 5 (and anything
# below) down by one. This turns the old row 5 into row 6, matching the
# diff's duplication of the original row 5 data into row 6.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the updated weekly price entry.
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(5, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(5, 4).Value = 44533
$ws.Cells.Item(5, 5).Value = 15
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100103
$ws.Cells.Item(5, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(5, 9).Value = 100103003
$ws.Cells.Item(5, 10).Value = "Damasco"
$ws.Cells.Item(5, 11).Value = "Castle Brite"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 140
$ws.Cells.Item(5, 14).Value = 14000
$ws.Cells.Item(5, 15).Value = 15000
$ws.Cells.Item(5, 16).Value = 14500
$ws.Cells.Item(5, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(5, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(5, 19).Value = 1450
$ws.Cells.Item(5, 20).Value = 10
